$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to a literal text value even when the text looks numeric,
# then restore the default (unstyled) cell style so no formatting side effects remain.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "28.699.27"
$ws.Range("E2").Value = "  +1.39%  "
$ws.Range("D3").Value = "1.870.21"
$ws.Range("E3").Value = "  +1.64%  "
Set-TextValue $ws.Range("D4") "1.005"
$ws.Range("E4").Value = "  +0.23%  "
Set-TextValue $ws.Range("D5") "326.84"
$ws.Range("E5").Value = "  -1.59%  "
$ws.Range("E6").Value = "  +0.32%  "
Set-TextValue $ws.Range("D7") "0.4651"
$ws.Range("E7").Value = "  +0.60%  "
Set-TextValue $ws.Range("D8") "0.3919"
$ws.Range("E8").Value = "  +1.35%  "
Set-TextValue $ws.Range("D9") "0.07916"
$ws.Range("E9").Value = "  +0.69%  "
Set-TextValue $ws.Range("D10") "0.9734"
$ws.Range("E10").Value = "  +0.44%  "
Set-TextValue $ws.Range("D11") "22.34"
$ws.Range("E11").Value = "  +1.72%  "
$ws.Range("D12").Value = "1.887.66"
$ws.Range("E12").Value = "  +0.25%  "
Set-TextValue $ws.Range("D13") "5.732"
$ws.Range("E13").Value = "  -0.52%  "
Set-TextValue $ws.Range("D14") "6.945"
$ws.Range("E14").Value = "  +0.16%  "
Set-TextValue $ws.Range("D15") "0.06903"
$ws.Range("E15").Value = "  +0.37%  "
Set-TextValue $ws.Range("D16") "88.79"
$ws.Range("E16").Value = "  +1.84%  "
$ws.Range("E17").Value = "  +0.35%  "
Set-TextValue $ws.Range("D18") "0.00001004"
$ws.Range("E18").Value = "  +0.80%  "
Set-TextValue $ws.Range("D19") "16.96"
$ws.Range("E19").Value = "  +0.27%  "
Set-TextValue $ws.Range("D20") "1.005"
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("D21").Value = "28.717.53"
$ws.Range("E21").Value = "  +1.36%  "
Set-TextValue $ws.Range("D22") "5.330"
$ws.Range("E22").Value = "  -0.52%  "
Set-TextValue $ws.Range("D23") "11.09"
$ws.Range("E23").Value = "  -0.13%  "
Set-TextValue $ws.Range("D24") "2.129"
$ws.Range("E24").Value = "  -1.59%  "
$ws.Range("D25").Value = "2.159.88"
$ws.Range("E25").Value = "  +4.38%  "
Set-TextValue $ws.Range("D26") "155.23"
$ws.Range("E26").Value = "  +1.03%  "
Set-TextValue $ws.Range("D27") "19.29"
$ws.Range("E27").Value = "  +0.17%  "
Set-TextValue $ws.Range("D28") "5.760"
$ws.Range("E28").Value = "  -2.02%  "
Set-TextValue $ws.Range("D29") "1.993"
$ws.Range("E29").Value = "  +0.75%  "
Set-TextValue $ws.Range("D30") "119.41"
$ws.Range("E30").Value = "  +2.17%  "
Set-TextValue $ws.Range("D31") "0.09344"
$ws.Range("E31").Value = "  +0.06%  "
Set-TextValue $ws.Range("D32") "0.9409"
$ws.Range("E32").Value = "  -0.67%  "
Set-TextValue $ws.Range("D33") "5.328"
$ws.Range("E33").Value = "  +0.41%  "
Set-TextValue $ws.Range("D34") "1.344"
$ws.Range("E34").Value = "  +0.87%  "
Set-TextValue $ws.Range("D35") "3.353"
$ws.Range("E35").Value = "  -2.88%  "
Set-TextValue $ws.Range("D36") "0.05830"
$ws.Range("E36").Value = "  -3.84%  "
Set-TextValue $ws.Range("D37") "0.02114"
$ws.Range("E37").Value = "  -2.73%  "
Set-TextValue $ws.Range("D38") "1.156"
$ws.Range("E38").Value = "  -0.37%  "
Set-TextValue $ws.Range("D39") "7.953"
$ws.Range("E39").Value = "  +4.84%  "
Set-TextValue $ws.Range("D40") "0.5669"
$ws.Range("E40").Value = "  +0.60%  "
Set-TextValue $ws.Range("D41") "9.970"
$ws.Range("E41").Value = "  -0.79%  "
Set-TextValue $ws.Range("D42") "0.1776"
$ws.Range("E42").Value = "  -0.63%  "
Set-TextValue $ws.Range("D43") "0.07313"
$ws.Range("E43").Value = "  +3.90%  "
Set-TextValue $ws.Range("D44") "2.250"
$ws.Range("E44").Value = "  -4.83%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws.Range("D45") "0.5339"
$ws.Range("E45").Value = "  +0.32%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D46") "11.70"
$ws.Range("E46").Value = "  +0.31%  "
$ws.Range("E47").Value = "  -6.18%  "
Set-TextValue $ws.Range("D48") "1.851"
$ws.Range("E48").Value = "  -0.04%  "
Set-TextValue $ws.Range("D49") "113.83"
$ws.Range("E49").Value = "  +0.29%  "
Set-TextValue $ws.Range("D50") "2.357"
$ws.Range("E50").Value = "  +1.47%  "
$ws.Range("E51").Value = "  +0.40%  "
